# Add a bold heading row to the users sheet, and append a new user record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top for column headings; this pushes all
# existing data down by one row (old row N -> new row N+1).
$ws.Rows.Item(1).Insert()

$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Email"
$ws.Range("C1").Value = "Phone"
$ws.Range("D1").Value = "Address"
$ws.Range("E1").Value = "Course"
$ws.Range("A1:E1").Font.Bold = $true

# Append the new user as the last row of the table.
$ws.Range("A8").Value = "Sseguya John"
$ws.Range("C8").Value = 256772554007
$ws.Range("D8").Value = "Nansana"
$ws.Range("E8").Value = "Certificate in Metal Works"

$ws.Range("E8").Select() | Out-Null
